$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2 through 140). All of these currently hold 45175 and need to be
# updated to 45177, keeping their existing number formatting/style.
for ($row = 2; $row -le 140; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45175) {
        $cell.Value2 = 45177
    }
}
